$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "subtylostyle"
$ws.Range("A3").Value = "thin spiraster"
